$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "scale_factors" row (99) and the "energies" row (100) in the
# AdvancedRadiograph field table were in the wrong order - swap them
# (values + formatting) so "energies" comes first, matching the rest of
# the spec. Using Range.Sort (rather than copying values cell by cell)
# keeps each cell's original formatting attached as the rows trade places.
$swapRange = $ws.Range("B99:E100")
$sortKey = $ws.Range("B99:B100")
$swapRange.Sort($sortKey, 1)

# Restore the sheet's on-screen scroll position / selection to where the
# author left it after making the edit.
$ws.Activate()
$ws.Range("C101").Select()
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
